$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 61 (G61=4604)
$ws.Range("H61").Value = 250
$ws.Range("I61").Value = 250
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 750
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -578
$ws.Range("N61").ClearContents()

# Row 100 (G100=19906)
$ws.Range("H100").Value = 2529.3845
$ws.Range("I100").Value = 2273.2856
$ws.Range("J100").Value = 2828.1667
$ws.Range("K100").Value = 2273.2856
$ws.Range("L100").Value = 2828.1667
$ws.Range("M100").Value = -1732.2856
$ws.Range("N100").Value = -3910.1667

# Row 137 (G137=44013)
$ws.Range("H137").Value = 3066
$ws.Range("J137").Value = 3385.5
$ws.Range("L137").Value = 10156.5
$ws.Range("N137").Value = -15256.5

$ws = $wb.Worksheets.Item("ARM")
# Row 34 (G34=2753)
$ws.Range("H34").Value = 49999.5
$ws.Range("I34").Value = 49999.5
$ws.Range("K34").Value = 49999.5
$ws.Range("M34").Value = -49728.5

# Row 62 (G62=10719)
$ws.Range("H62").Value = 90000
$ws.Range("J62").Value = 90000
$ws.Range("L62").Value = 90000
$ws.Range("N62").Value = -91248

# Row 63 (G63=12528)
$ws.Range("H63").Value = 1780.4375
$ws.Range("I63").Value = 1238.3334
$ws.Range("J63").Value = 1905.5385
$ws.Range("K63").Value = 1238.3334
$ws.Range("L63").Value = 1905.5385
$ws.Range("M63").Value = -552.3334
$ws.Range("N63").Value = -3277.5385

# Row 65 (G65=10719)
$ws.Range("H65").Value = 90000
$ws.Range("J65").Value = 90000
$ws.Range("L65").Value = 270000
$ws.Range("N65").Value = -276240

# Row 66 (G66=12528)
$ws.Range("H66").Value = 1780.4375
$ws.Range("I66").Value = 1238.3334
$ws.Range("J66").Value = 1905.5385
$ws.Range("K66").Value = 6191.666999999999
$ws.Range("L66").Value = 9527.692500000001
$ws.Range("M66").Value = -2759.666999999999
$ws.Range("N66").Value = -16391.6925

# Row 110 (G110=27708)
$ws.Range("H110").Value = 4478.6113
$ws.Range("I110").Value = 3897.9333
$ws.Range("K110").Value = 3897.9333
$ws.Range("M110").Value = -1852.9333

# Row 122 (G122=36168)
$ws.Range("H122").Value = 1405.5555
$ws.Range("I122").Value = 1400.625
$ws.Range("K122").Value = 4201.875
$ws.Range("M122").Value = -1751.875

# Row 132 (G132=43997)
$ws.Range("H132").Value = 25003436
$ws.Range("I132").Value = 2607.3125
$ws.Range("K132").Value = 7821.9375
$ws.Range("M132").Value = -5291.9375

$ws = $wb.Worksheets.Item("BSM")
# Row 82 (G82=11877)
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

# Row 85 (G85=11877)
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

# Row 86 (G86=12526)
$ws.Range("H86").Value = 2345.4
$ws.Range("I86").Value = 2231.75
$ws.Range("K86").Value = 2231.75
$ws.Range("M86").Value = -1108.75

# Row 88 (G88=10626)
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# Row 89 (G89=12526)
$ws.Range("H89").Value = 2345.4
$ws.Range("I89").Value = 2231.75
$ws.Range("K89").Value = 11158.75
$ws.Range("M89").Value = -5542.75

# Row 91 (G91=10626)
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# Row 99 (G99=19943)
$ws.Range("H99").Value = 1953.4445
$ws.Range("I99").Value = 1997
$ws.Range("J99").Value = 1948
$ws.Range("K99").Value = 1997
$ws.Range("L99").Value = 1948
$ws.Range("M99").Value = -499
$ws.Range("N99").Value = -4944

# Row 105 (G105=19947)
$ws.Range("H105").Value = 3384.9
$ws.Range("I105").Value = 3756
$ws.Range("K105").Value = 3756
$ws.Range("M105").Value = -2009

# Row 107 (G107=27706)
$ws.Range("H107").Value = 54276.5
$ws.Range("I107").Value = 16415.25
$ws.Range("J107").Value = 129999
$ws.Range("K107").Value = 16415.25
$ws.Range("L107").Value = 129999
$ws.Range("M107").Value = -14495.25
$ws.Range("N107").Value = -133839

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G31=44023)
$ws.Range("H31").Value = 1849.2222
$ws.Range("I31").Value = 1081
$ws.Range("J31").Value = 2068.7144
$ws.Range("K31").Value = 1081
$ws.Range("L31").Value = 2068.7144
$ws.Range("M31").Value = -786
$ws.Range("N31").Value = -2658.7144

# Row 34 (G34=44023)
$ws.Range("H34").Value = 1849.2222
$ws.Range("I34").Value = 1081
$ws.Range("J34").Value = 2068.7144
$ws.Range("K34").Value = 1081
$ws.Range("L34").Value = 2068.7144
$ws.Range("M34").Value = -879
$ws.Range("N34").Value = -2472.7144

# Row 86 (G86=12584)
$ws.Range("H86").Value = 4199.25
$ws.Range("J86").Value = 4499.5
$ws.Range("L86").Value = 4499.5
$ws.Range("N86").Value = -6745.5

# Row 89 (G89=12584)
$ws.Range("H89").Value = 4199.25
$ws.Range("J89").Value = 4499.5
$ws.Range("L89").Value = 22497.5
$ws.Range("N89").Value = -33729.5

# Row 99 (G99=36198)
$ws.Range("H99").Value = 1182.875
$ws.Range("I99").Value = 914.3333
$ws.Range("K99").Value = 914.3333
$ws.Range("M99").Value = 583.6667

# Row 122 (G122=36196)
$ws.Range("H122").Value = 17822.666
$ws.Range("I122").Value = 804.7778
$ws.Range("J122").Value = 43349.5
$ws.Range("K122").Value = 2414.3334
$ws.Range("L122").Value = 130048.5
$ws.Range("M122").Value = 35.66660000000002
$ws.Range("N122").Value = -134948.5

# Row 126 (G126=36198)
$ws.Range("H126").Value = 1182.875
$ws.Range("I126").Value = 914.3333
$ws.Range("K126").Value = 2742.9999
$ws.Range("M126").Value = -272.9998999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 64 (G64=12861)
$ws.Range("H64").Value = 25007
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 25007
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 75021
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -75561

# Row 67 (G67=12861)
$ws.Range("H67").Value = 25007
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 25007
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 75021
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -76893

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (G80=12521)
$ws.Range("H80").Value = 2995.4285
$ws.Range("I80").Value = 2259.6667
$ws.Range("K80").Value = 2259.6667
$ws.Range("M80").Value = -1261.6667

# Row 83 (G83=12521)
$ws.Range("H83").Value = 2995.4285
$ws.Range("I83").Value = 2259.6667
$ws.Range("K83").Value = 11298.3335
$ws.Range("M83").Value = -6306.333500000001

# Row 122 (G122=36182)
$ws.Range("H122").Value = 4014.2144
$ws.Range("I122").Value = 4638.4
$ws.Range("J122").Value = 2453.75
$ws.Range("K122").Value = 13915.2
$ws.Range("L122").Value = 7361.25
$ws.Range("M122").Value = -11465.2
$ws.Range("N122").Value = -12261.25

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (G132=44058)
$ws.Range("H132").Value = 4966.3335
$ws.Range("I132").Value = 3452
$ws.Range("J132").Value = 7995
$ws.Range("K132").Value = 10356
$ws.Range("L132").Value = 23985
$ws.Range("M132").Value = -7826
$ws.Range("N132").Value = -29045

$ws = $wb.Worksheets.Item("WVR")
# Row 58 (G58=3187)
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
